$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header (F1) -- triggers new shared string "cell_volume_on_filter_uL"
$ws.Range("F1").Value = "cell_volume_on_filter_uL"

# Column F width
$ws.Columns("F").ColumnWidth = 23.67

# Column F cell values (per-organism measurements)
$ws.Range("F37").Value = 0.11362
$ws.Range("F38").Value = 0.10396
$ws.Range("F39").Value = 0.09982
$ws.Range("F44").Value = 0.31850000000000001
$ws.Range("F45").Value = 0.34858
$ws.Range("F46").Value = 0.31701000000000001
$ws.Range("F47").Value = 2.20472
$ws.Range("F48").Value = 2.1892800000000001
$ws.Range("F49").Value = 0.70799999999999996
$ws.Range("F50").Value = 0.75483999999999996
$ws.Range("F51").Value = 3.7112579999999999
$ws.Range("F52").Value = 5.3510318180000001
$ws.Range("F53").Value = 3.5155285709999999
$ws.Range("F54").Value = 0.45255630899999999
$ws.Range("F55").Value = 0.48748622200000002
$ws.Range("F56").Value = 0.46708456300000001
$ws.Range("F57").Value = 0.39088334299999999
$ws.Range("F58").Value = 0.34649945900000001
$ws.Range("F59").Value = 0.39507073599999998
$ws.Range("F60").Value = 0.062888292
$ws.Range("F61").Value = 0.054405
$ws.Range("F62").Value = 0.054925488
$ws.Range("F63").Value = 7.0207955999999996
$ws.Range("F64").Value = 2.2855254550000002
$ws.Range("F65").Value = 4.9629761999999999
$ws.Range("F69").Value = 4.4295660000000003
$ws.Range("F70").Value = 5.0769299999999999
$ws.Range("F71").Value = 5.6976899999999997
$ws.Range("F72").Value = 3.9806896549999999
$ws.Range("F73").Value = 4.1544827590000004
$ws.Range("F74").Value = 2.5326315789999998
$ws.Range("F75").Value = 52.498559999999998
$ws.Range("F76").Value = 55.513680000000001
$ws.Range("F77").Value = 10.10952
$ws.Range("F78").Value = 13.268750000000001
$ws.Range("F79").Value = 12.20459
$ws.Range("F80").Value = 4.4022199999999998
$ws.Range("F81").Value = 3.1392699999999998
$ws.Range("F82").Value = 6.3636799999999996
$ws.Range("F83").Value = 0.33903
$ws.Range("F84").Value = 0.48124
$ws.Range("F85").Value = 0.55139000000000005
$ws.Range("F86").Value = 0.56953487999999997
$ws.Range("F87").Value = 0.58277279999999998
$ws.Range("F88").Value = 0.66859104000000003
$ws.Range("F89").Value = 0.557926221
$ws.Range("F90").Value = 0.57103792399999997
$ws.Range("F91").Value = 0.53643271299999995
$ws.Range("F92").Value = 0.45029728000000002
$ws.Range("F93").Value = 0.39780863999999999
$ws.Range("F94").Value = 0.23205503999999999
$ws.Range("F95").Value = 2.75917726
$ws.Range("F96").Value = 5.1155549200000001
$ws.Range("F97").Value = 1.5322825600000001
$ws.Range("F98").Value = 0.169351485
$ws.Range("F99").Value = 0.17240298500000001
$ws.Range("F100").Value = 0.175519802
$ws.Range("F104").Value = 0.36580499999999999
$ws.Range("F105").Value = 0.38092090899999997
$ws.Range("F106").Value = 0.29425636399999999
$ws.Range("F107").Value = 0.29131000000000001
$ws.Range("F108").Value = 0.33255000000000001
$ws.Range("F109").Value = 0.381324
$ws.Range("F110").Value = 2.6808257310632899
$ws.Range("F111").Value = 2.6808257310632899
$ws.Range("F112").Value = 2.6808257310632899
$ws.Range("F113").Value = 0.65449846949787349
$ws.Range("F114").Value = 0.65449846949787349
$ws.Range("F115").Value = 0.65449846949787349
$ws.Range("F116").Value = 0.33510321638291124
$ws.Range("F117").Value = 0.33510321638291124
$ws.Range("F118").Value = 0.33510321638291124

# Highlight fill + number format for newly back-filled D/E cells
$ws.Range("D66").Interior.Color = 65535
$ws.Range("E66").Interior.Color = 65535
$ws.Range("D67").Interior.Color = 65535
$ws.Range("E67").Interior.Color = 65535
$ws.Range("D68").Interior.Color = 65535
$ws.Range("E68").Interior.Color = 65535
$ws.Range("D101").Interior.Color = 65535
$ws.Range("E101").Interior.Color = 65535
$ws.Range("E101").NumberFormat = "0.00E+00"
$ws.Range("E101").Value = 10000000
$ws.Range("D102").Interior.Color = 65535
$ws.Range("E102").Interior.Color = 65535
$ws.Range("E102").NumberFormat = "0.00E+00"
$ws.Range("E102").Value = 10000000
$ws.Range("D103").Interior.Color = 65535
$ws.Range("E103").Interior.Color = 65535
$ws.Range("E103").NumberFormat = "0.00E+00"
$ws.Range("E103").Value = 10000000
$ws.Range("D110").Interior.Color = 65535
$ws.Range("E110").Interior.Color = 65535
$ws.Range("E110").NumberFormat = "0.00E+00"
$ws.Range("E110").Value = 10000000
$ws.Range("D111").Interior.Color = 65535
$ws.Range("E111").Interior.Color = 65535
$ws.Range("E111").NumberFormat = "0.00E+00"
$ws.Range("E111").Value = 10000000
$ws.Range("D112").Interior.Color = 65535
$ws.Range("E112").Interior.Color = 65535
$ws.Range("E112").NumberFormat = "0.00E+00"
$ws.Range("E112").Value = 10000000
$ws.Range("D113").Interior.Color = 65535
$ws.Range("E113").Interior.Color = 65535
$ws.Range("E113").NumberFormat = "0.00E+00"
$ws.Range("E113").Value = 10000000
$ws.Range("D114").Interior.Color = 65535
$ws.Range("E114").Interior.Color = 65535
$ws.Range("E114").NumberFormat = "0.00E+00"
$ws.Range("E114").Value = 10000000
$ws.Range("D115").Interior.Color = 65535
$ws.Range("E115").Interior.Color = 65535
$ws.Range("E115").NumberFormat = "0.00E+00"
$ws.Range("E115").Value = 10000000
$ws.Range("D116").Interior.Color = 65535
$ws.Range("E116").Interior.Color = 65535
$ws.Range("E116").NumberFormat = "0.00E+00"
$ws.Range("E116").Value = 10000000
$ws.Range("D117").Interior.Color = 65535
$ws.Range("E117").Interior.Color = 65535
$ws.Range("E117").NumberFormat = "0.00E+00"
$ws.Range("E117").Value = 10000000
$ws.Range("D118").Interior.Color = 65535
$ws.Range("E118").Interior.Color = 65535
$ws.Range("E118").NumberFormat = "0.00E+00"
$ws.Range("E118").Value = 10000000
$ws.Range("D119").Interior.Color = 65535
$ws.Range("E119").Interior.Color = 65535
$ws.Range("E119").NumberFormat = "0.00E+00"
$ws.Range("E119").Value = 10000000
$ws.Range("D120").Interior.Color = 65535
$ws.Range("E120").Interior.Color = 65535
$ws.Range("E120").NumberFormat = "0.00E+00"
$ws.Range("E120").Value = 10000000
$ws.Range("D121").Interior.Color = 65535
$ws.Range("E121").Interior.Color = 65535
$ws.Range("E121").NumberFormat = "0.00E+00"
$ws.Range("E121").Value = 10000000
$ws.Range("D122").Interior.Color = 65535
$ws.Range("E122").Interior.Color = 65535
$ws.Range("E122").NumberFormat = "0.00E+00"
$ws.Range("E122").Value = 10000000
$ws.Range("D123").Interior.Color = 65535
$ws.Range("E123").Interior.Color = 65535
$ws.Range("E123").NumberFormat = "0.00E+00"
$ws.Range("E123").Value = 10000000
$ws.Range("D124").Interior.Color = 65535
$ws.Range("E124").Interior.Color = 65535
$ws.Range("E124").NumberFormat = "0.00E+00"
$ws.Range("E124").Value = 10000000

# View state (best effort; selection reflects final cursor position from the source edit)
$ws.Range("F121").Select()
